# Applies the updated coinranking.com snapshot values (price/volume refresh,
# plus two coin-pair row swaps) to the crypto listing on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.311.46'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.18%  '

# Row 3: Ethereum
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.685.49'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.90%  '

# Row 4: TetherUSD
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.41%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.27'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.44%  '

# Row 7: USDC
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.04%  '

# Row 8: XRP
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.50%  '

# Row 9: Dogecoin
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +10.14%  '

# Row 10: Toncoin
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.01'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.40%  '

# Row 11: Cardano
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.409'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.27%  '

# Row 12: TRON
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.70%  '

# Row 13: ShibaInu
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000214'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +24.03%  '

# Row 14: Avalanche
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '30.72'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.30%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.169.58'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.07%  '

# Row 16: WrappedBTC
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.094.77'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.08%  '

# Row 17: WrappedEther
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.681.35'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.29%  '

# Row 18: Chainlink
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.74'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.87%  '

# Row 19: Polkadot
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.73%  '

# Row 20: BitcoinCash
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '362.74'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.50%  '

# Row 21: Uniswap
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.27%  '

# Row 22: Dai
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.26%  '

# Row 23: Litecoin
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.30'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.73%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.76'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.00%  '

# Row 25: PEPE
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +17.23%  '

# Row 26: SuiNetwork
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.68'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.98%  '

# Row 27: Kaspa
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.175'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.82%  '

# Row 28: Fetch.AI
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.67'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.90%  '

# Row 29: Aptos
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.22'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.22%  '

# Row 30: PancakeSwap
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +7.55%  '

# Row 31: Bittensor
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '542.30'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.27%  '

# Row 32: Binance-PegBSC-USD
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.21%  '

# Row 33: ImmutableX
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.43%  '

# Row 34: RenderToken
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.97%  '

# Row 35: NEARProtocol
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.58'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.58%  '

# Row 36: PolygonEcosystemToken
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.59%  '

# Row 37: EthereumClassic
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.86'

# Row 38: Stacks
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.36%  '

# Row 39: Monero
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '163.10'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.99%  '

# Row 40: FirstDigitalUSD
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.01%  '

# Row 41: Aave
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '170.98'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.49%  '

# Row 42: USDe
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.08%  '

# Row 43: OKB
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.53'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.54%  '

# Row 44: Filecoin
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.27%  '

# Row 45: dogwifhat
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.36'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.89%  '

# Row 46: Hedera
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.00%  '

# Row 47: InjectiveProtocol
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.46%  '

# Row 48: Mantle
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.664'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.20%  '

# Row 49: VeChain
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.38%  '

# Row 50: EnergySwap
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.07'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.40%  '

# Row 51: Stellar
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0992'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.81%  '
